# Update the "Förändrad" (Changed) date column (C) for all data rows.
# The value 45803 (2025-05-26) is incremented by one day to 45804 (2025-05-27)
# for every row in the table (rows 2 through 43 on the active sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45803) {
        $cell.Value2 = 45804
    }
}
